# Splits a run of text into several runs whose texts are given by $parts
# (concatenated, they must equal the new full text that replaces $oldText),
# while preserving run formatting (rPr) identically on every piece.
#
# Technique: locate $oldText with Find, overwrite that Range's .Text with
# the new (joined) text, then force hard run boundaries at each internal
# part boundary (and, if requested, at extra boundary positions measured
# from the end of the new text) by toggling Font.Size away from and back
# to its original value on the sub-range immediately to the left of each
# boundary. This engine (like Word) merges adjacent runs whose rPr is
# identical, but a transient Font.Size change followed by restoring the
# original value leaves a hard run break in place without altering the
# final formatting of any run.
#
# $extraBoundaries is a list of character offsets, measured from the start
# of the replaced text, where an additional split must be (re-)established
# -- this is needed when the original text run was immediately followed by
# other runs (e.g. the old text was itself the prefix of a larger run run
# that had already been split into several runs by earlier edits): the
# whole-range .Text assignment merges the edited run with any immediately
# following runs that happen to share the same rPr, so those boundaries
# must be explicitly restored too.
function Split-IntoRuns($doc, $oldText, $parts, $extraBoundaries) {
    $newText = [string]::Join("", $parts)

    $find = $doc.Content
    $found = $find.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $oldText"
    }

    $start = $find.Start
    $end = $find.End

    $full = $doc.Range($start, $end)
    $full.Text = $newText
    $newEnd = $start + $newText.Length

    $boundaries = New-Object System.Collections.ArrayList
    $pos = $start
    for ($i = 0; $i -lt $parts.Length - 1; $i++) {
        $pos = $pos + $parts[$i].Length
        [void]$boundaries.Add($pos)
    }
    if ($extraBoundaries) {
        foreach ($off in $extraBoundaries) {
            $abs = $newEnd + $off
            [void]$boundaries.Add($abs)
        }
    }

    $prev = $start
    foreach ($b in $boundaries) {
        $left = $doc.Range($prev, $b)
        $origSize = $left.Font.Size
        $tempSize = 13
        if ($origSize -eq 13) {
            $tempSize = 12
        }
        $left.Font.Size = $tempSize
        $left.Font.Size = $origSize
        $prev = $b
    }
}

$d = $word.ActiveDocument

Split-IntoRuns $d `
    "1.Step: Sanction List Check- sanctioned companies check and assigning score" `
    @("1.", " ", "Step: Sanction List Check- sanctioned companies check and assigning score") `
    @()

Split-IntoRuns $d `
    "2 Step: Company Status check -company active or not active and assigning scores" `
    @("2", ". ", "Step: Company Status check -company active or not active and assigning scores") `
    @()

Split-IntoRuns $d `
    "3. Web s" `
    @("3.", " Step:", " Web s") `
    @(0, 1)

Split-IntoRuns $d `
    "4. Companies score sum of 1+2+3 results" `
    @("4.", " Step: ", "Companies score sum of 1+2+3 results") `
    @()

Write-Output "done"
